# Auto-generated edit script
# Commit: Add data for 2025-09-06
# Updates partial-year 2025 (column L) violent crime counts across citywide,
# by-neighborhood summary, and individual neighborhood detail sheets, reflecting
# newly recorded incidents as of 2025-09-06 (plus a couple of 2022 total corrections).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 4593
$ws.Range("L3").Value = 4928
$ws.Range("I4").Value = 1844
$ws.Range("L4").Value = 1221
$ws.Range("L5").Value = 285
$ws.Range("L6").Value = 4211
$ws.Range("I7").Value = 26314
$ws.Range("L7").Value = 15238

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("L3").Value = 10
$ws.Range("L7").Value = 39

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L3").Value = 340
$ws.Range("L4").Value = 76
$ws.Range("L6").Value = 268

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 189
$ws.Range("L3").Value = 235
$ws.Range("L4").Value = 41
$ws.Range("L6").Value = 215
$ws.Range("L7").Value = 694

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L2").Value = 78
$ws.Range("L7").Value = 206

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 167
$ws.Range("L3").Value = 195
$ws.Range("L6").Value = 152
$ws.Range("L7").Value = 562

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L3").Value = 93
$ws.Range("L6").Value = 80
$ws.Range("L7").Value = 296

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L4").Value = 21
$ws.Range("L7").Value = 261

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("L6").Value = 30
$ws.Range("L7").Value = 71

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 127
$ws.Range("L6").Value = 117
$ws.Range("L16").Value = 32
$ws.Range("L17").Value = 28
$ws.Range("L18").Value = 107
$ws.Range("L19").Value = 419
$ws.Range("L20").Value = 388
$ws.Range("L23").Value = 165
$ws.Range("L27").Value = 135
$ws.Range("L29").Value = 844
$ws.Range("L30").Value = 71
$ws.Range("L33").Value = 694
$ws.Range("L36").Value = 199
$ws.Range("L37").Value = 562
$ws.Range("L47").Value = 107
$ws.Range("L48").Value = 197
$ws.Range("L49").Value = 79
$ws.Range("L51").Value = 187
$ws.Range("L54").Value = 319
$ws.Range("L57").Value = 54
$ws.Range("I63").Value = 265
$ws.Range("L63").Value = 45
$ws.Range("L64").Value = 105
$ws.Range("L65").Value = 296
$ws.Range("L67").Value = 525
$ws.Range("L68").Value = 49
$ws.Range("L69").Value = 39
$ws.Range("L71").Value = 43
$ws.Range("L78").Value = 200
$ws.Range("L84").Value = 147
$ws.Range("L85").Value = 781
$ws.Range("L86").Value = 115
$ws.Range("L89").Value = 218
$ws.Range("L90").Value = 153
$ws.Range("L94").Value = 191
$ws.Range("L95").Value = 206
$ws.Range("L96").Value = 172
$ws.Range("L99").Value = 261
$ws.Range("I101").Value = 26314
$ws.Range("L101").Value = 15238

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 152
$ws.Range("L3").Value = 204
$ws.Range("L6").Value = 120
$ws.Range("L7").Value = 525

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L3").Value = 53
$ws.Range("L7").Value = 147

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("L6").Value = 32
$ws.Range("L7").Value = 79

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L6").Value = 155
$ws.Range("L7").Value = 319

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L3").Value = 318
$ws.Range("L7").Value = 844

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L3").Value = 48
$ws.Range("L7").Value = 197

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 148
$ws.Range("L3").Value = 129
$ws.Range("L6").Value = 121
$ws.Range("L7").Value = 419

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L3").Value = 33
$ws.Range("L5").Value = 1
$ws.Range("L7").Value = 117

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 141
$ws.Range("L6").Value = 137

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L6").Value = 61
$ws.Range("L7").Value = 200

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L3").Value = 60
$ws.Range("L7").Value = 165

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L6").Value = 49
$ws.Range("L7").Value = 172

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L5").Value = 3
$ws.Range("L7").Value = 105

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L3").Value = 127
$ws.Range("L7").Value = 388

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("L2").Value = 39
$ws.Range("L3").Value = 38
$ws.Range("L6").Value = 17
$ws.Range("L7").Value = 107

$ws = $wb.Worksheets.Item("Burnside")
$ws.Range("L3").Value = 10
$ws.Range("L7").Value = 28

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L3").Value = 58
$ws.Range("L6").Value = 51
$ws.Range("L7").Value = 199

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L2").Value = 43
$ws.Range("L7").Value = 191

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L2").Value = 40
$ws.Range("L7").Value = 107

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L3").Value = 42
$ws.Range("L6").Value = 30
$ws.Range("L7").Value = 127

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L3").Value = 64
$ws.Range("L7").Value = 218

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L3").Value = 41
$ws.Range("L6").Value = 40
$ws.Range("L7").Value = 135

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("L4").Value = 63
$ws.Range("L7").Value = 115

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L2").Value = 51
$ws.Range("L3").Value = 46
$ws.Range("L7").Value = 153

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L3").Value = 61
$ws.Range("L6").Value = 42
$ws.Range("L7").Value = 187

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("L2").Value = 16
$ws.Range("L7").Value = 49

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("L6").Value = 16
$ws.Range("L7").Value = 54

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 236
$ws.Range("L5").Value = 18
$ws.Range("L6").Value = 165
$ws.Range("L7").Value = 781

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("L6").Value = 10
$ws.Range("L7").Value = 43

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("L6").Value = 21
$ws.Range("L7").Value = 32
